$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new analysis columns (G, H, I) with headers, continuing the
# existing shared-string table.
$ws.Range("G1").Value = "Average Error (lbf)"
$ws.Range("H1").Value = "Max Error (lbf)"
$ws.Range("I1").Value = "Min Error (lbf)"

# Reuse the exact header formatting already used by A1:E1 (bold, centered
# horizontally/vertically) by copying its format onto the new headers.
$ws.Range("A1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

# Best-fit column widths for the new columns (character widths chosen so the
# stored sheet width matches the target: 15, 12, and ~11.83).
$ws.Columns.Item(7).ColumnWidth = 15
$ws.Columns.Item(8).ColumnWidth = 12
$ws.Columns.Item(9).ColumnWidth = 10.998697916666666

# Move/restore the active selection to match the author's final cursor spot.
$ws.Range("G16").Select()
